$wb = $excel.ActiveWorkbook

$msg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f82830a6871036ae19c4ee0d5c9245745ef3eca9/e2e/a631560f-dea5-4037-ad31-4c39239a9508.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b164b25d7e2b014f0bbbcd4b5b7647d28ea59cd8/e2e/a631560f-dea5-4037-ad31-4c39239a9508.md."

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "a631560f-dea5-4037-ad31-4c39239a9508.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f82830a6871036ae19c4ee0d5c9245745ef3eca9/e2e/a631560f-dea5-4037-ad31-4c39239a9508.md", "", "", "a631560f-dea5-4037-ad31-4c39239a9508.md") | Out-Null

$wsZh.Range("J7").Value = "a631560f-dea5-4037-ad31-4c39239a9508.d26cca369428f3279109f05fe9c52e349a50243b.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-06 15:52:16"
$wsZh.Range("P7").Value = $msg

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "a631560f-dea5-4037-ad31-4c39239a9508.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f82830a6871036ae19c4ee0d5c9245745ef3eca9/e2e/a631560f-dea5-4037-ad31-4c39239a9508.md", "", "", "a631560f-dea5-4037-ad31-4c39239a9508.md") | Out-Null

$wsDe.Range("J7").Value = "a631560f-dea5-4037-ad31-4c39239a9508.d26cca369428f3279109f05fe9c52e349a50243b.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-06 15:52:28"
$wsDe.Range("P7").Value = $msg

Write-Host "Report generated for handback."
